# Adds two bullet points at the end of the "Additional features added" list:
#   1. Fills the existing (empty) trailing list paragraph with the
#      "green cards" bullet text.
#   2. Inserts a brand new list paragraph right after it for the
#      "purple cards" bullet text.
#
# Navigation is anchor-based (via Find + Next()) rather than hard-coded
# paragraph indices, so it is robust to the exact paragraph numbering.

$d = $word.ActiveDocument

# Locate the last populated bullet in the list ("...selectable in GUI").
# The paragraph immediately after it is the empty list-style paragraph
# that needs the new "green cards" text.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Added indicator if card is selectable in GUI",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph 'Added indicator if card is selectable in GUI'"
}

$anchorPara = $rng.Paragraphs(1)
$greenPara = $anchorPara.Next()

# Fill the existing empty list paragraph with the green-card bullet text.
$greenPara.Range.InsertAfter("Added card effect for green cards")

# Insert a new list paragraph right after it (inherits the same
# ListParagraph style / numbering) and fill it with the purple-card text.
$greenPara.Range.InsertParagraphAfter()
$purplePara = $greenPara.Next()
$purplePara.Range.InsertAfter("Partially added card effects for purple cards (except for V points at end)")
